$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update shared string label in B6
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# Update numeric values for rows 4-13, columns C, D, E, G, H, I, J
$ws.Range("C4").Value = 0.1784339608490212
$ws.Range("D4").Value = 0.9600480012000301
$ws.Range("E4").Value = -0.7825365634140854
$ws.Range("G4").Value = -0.0637890947273682
$ws.Range("H4").Value = 0.4810620265506638
$ws.Range("I4").Value = 0.4464906622665568
$ws.Range("J4").Value = 0.8229850746268658

$ws.Range("C5").Value = 0.03470186754668867
$ws.Range("D5").Value = 0.09336233405835148
$ws.Range("E5").Value = 0.02301957548938724
$ws.Range("G5").Value = 0.08865071626790672
$ws.Range("H5").Value = 0.1521833045826146
$ws.Range("I5").Value = 0.134766369159229
$ws.Range("J5").Value = 0.235223880597015

$ws.Range("C6").Value = 0.07646441161029026
$ws.Range("D6").Value = 0.1813830345758644
$ws.Range("E6").Value = -0.02597464936623416
$ws.Range("G6").Value = 0.2306832670816771
$ws.Range("H6").Value = 0.1552118802970074
$ws.Range("I6").Value = 0.00109652741318533
$ws.Range("J6").Value = 0.04899872496812421

$ws.Range("C7").Value = 0.08621015525388136
$ws.Range("D7").Value = 0.08721218030450763
$ws.Range("E7").Value = 0.5809555238880973
$ws.Range("G7").Value = 0.4767209180229506
$ws.Range("H7").Value = 0.7291217280432011
$ws.Range("I7").Value = 0.8388689717242932
$ws.Range("J7").Value = -0.04908272706817671

$ws.Range("C8").Value = 0.83875196879922
$ws.Range("D8").Value = 0.04760519012975326
$ws.Range("E8").Value = 0.001354533863346584
$ws.Range("G8").Value = 0.01775144378609466
$ws.Range("H8").Value = 0.07909097727443187
$ws.Range("I8").Value = 0.06434410860271507
$ws.Range("J8").Value = 0.002715067876696918

$ws.Range("C9").Value = 0.3864666616665418
$ws.Range("D9").Value = -0.09345383634590865
$ws.Range("E9").Value = 0.03792994824870622
$ws.Range("G9").Value = -0.0276231905797645
$ws.Range("H9").Value = -0.05684392109802745
$ws.Range("I9").Value = -0.02090002250056252
$ws.Range("J9").Value = -0.06146103652591316

$ws.Range("C10").Value = 0.03397434935873397
$ws.Range("D10").Value = 0.01882547063676592
$ws.Range("E10").Value = -0.04248856221405535
$ws.Range("G10").Value = -0.0007560189004725118
$ws.Range("H10").Value = -0.01366834170854272
$ws.Range("I10").Value = -0.01665941648541214
$ws.Range("J10").Value = -0.005782644566114154

$ws.Range("C11").Value = -0.02046651166279157
$ws.Range("D11").Value = 0.08528163204080103
$ws.Range("E11").Value = 0.03320183004575115
$ws.Range("G11").Value = 0.08012300307507689
$ws.Range("H11").Value = 0.121752043801095
$ws.Range("I11").Value = 0.1309082727068177
$ws.Range("J11").Value = 0.01550588764719118

$ws.Range("C12").Value = 0.04957923948098703
$ws.Range("D12").Value = 0.07516387909697743
$ws.Range("E12").Value = 0.2638640966024151
$ws.Range("G12").Value = 0.7377919447986201
$ws.Range("H12").Value = 0.2309127728193205
$ws.Range("I12").Value = -0.0852126303157579
$ws.Range("J12").Value = 0.02504462611565289

$ws.Range("C13").Value = -0.1789529738243456
$ws.Range("D13").Value = 0.006000150003750095
$ws.Range("E13").Value = -0.01036375909397735
$ws.Range("G13").Value = -0.02665866646666167
$ws.Range("H13").Value = 0.03860496512412811
$ws.Range("I13").Value = 0.05545338633465837
$ws.Range("J13").Value = 0.005034125853146328
